# "fixed a text box" - slide 4, TextBox 13:
#  - remove the solid fill (tx2, lumMod 60%/lumOff 40%) so the box has no fill
#  - fix typo "transport" -> "transported" in the caption text

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(11)

# Correct the run's text without disturbing the trailing empty paragraph.
$run = $shape.TextFrame.TextRange.Runs(1, 1)
$run.Text = "Confused, Craig has entered a frenzy due to being transported here."

# Remove the solid background fill from the text box.
$shape.Fill.Visible = $false
